$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A, rows 2..97 currently hold labels "q1".."q96" (row R holds "q{R-1}").
# Decrement every label's numeric suffix by one: q1 -> q0, q2 -> q1, ..., q96 -> q95.
for ($row = 2; $row -le 97; $row++) {
    $newNumber = $row - 2
    $ws.Cells.Item($row, 1).Value = "q" + $newNumber
}
